$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Occasionally employed.global.pro"
$ws.Range("C1").Value = "Regularly employed.global.pro"
$ws.Range("D1").Value = "Student.global.pro"
$ws.Range("E1").Value = "Unemployed / discouraged.global.pro"
$ws.Range("F1").Value = "Receiving social benefits / pensioners / house-makers / disable.global.pro"
$ws.Range("G1").Value = "Other.global.pro"
$ws.Range("H1").Value = "Not known / missing.global.pro"
$ws.Range("I1").Value = "Total.global.pro"
